# Buttons für Farben hinzugefügt, Ändern der Farbe über Buttons implementiert,
# Variablennamen auf Englisch geändert.
#
# Concretely: mark the Java/Git identifiers in a few existing bullet points as
# "spell-check exceptions" (<w:proofErr w:type="spellStart"/>...<w:proofErr
# w:type="spellEnd"/> around their own run, split off from the surrounding
# German text), and append a new bullet point describing the new
# "change color via buttons" feature.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Pull the exact attributes (w14:paraId, w14:textId, w:rsidR, ...) currently
# on a paragraph's <w:p> element so the rewritten paragraph keeps them.
function Get-ParaAttrs($para) {
    $xml = $para.Range.WordOpenXML
    if ($xml -match '<w:p\s+([^>]*?)\s*/?>') {
        return $matches[1]
    }
    return ""
}

# Replace a paragraph's contents with new run/proofErr markup while keeping
# its existing list paragraph formatting (Listenabsatz / numId 1) and its
# original <w:p> attributes.
function Set-ListParaRuns($para, [string]$innerXml) {
    $attrs = Get-ParaAttrs $para
    $pPrXml = "<w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>"
    if ($attrs -ne "") {
        $openTag = "<w:p xmlns:w='$wNs' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' $attrs>"
    } else {
        $openTag = "<w:p xmlns:w='$wNs'>"
    }
    $xml = "$openTag$pPrXml$innerXml</w:p>"
    $para.Range.InsertXML($xml)
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -eq "Erstellen der GitRepository`r") {
        $inner = "<w:r><w:t xml:space='preserve'>Erstellen der </w:t></w:r>" +
                 "<w:proofErr w:type='spellStart'/>" +
                 "<w:r><w:t>GitRepository</w:t></w:r>" +
                 "<w:proofErr w:type='spellEnd'/>"
        Set-ListParaRuns $p $inner
    }
    elseif ($t -eq "Erstellen des JFrame`r") {
        $inner = "<w:r><w:t xml:space='preserve'>Erstellen des </w:t></w:r>" +
                 "<w:proofErr w:type='spellStart'/>" +
                 "<w:r><w:t>JFrame</w:t></w:r>" +
                 "<w:proofErr w:type='spellEnd'/>"
        Set-ListParaRuns $p $inner
    }
    elseif ($t -eq "Erstellen des JPanel`r") {
        $inner = "<w:r><w:t xml:space='preserve'>Erstellen des </w:t></w:r>" +
                 "<w:proofErr w:type='spellStart'/>" +
                 "<w:r><w:t>JPanel</w:t></w:r>" +
                 "<w:proofErr w:type='spellEnd'/>"
        Set-ListParaRuns $p $inner
    }
    elseif ($t -eq "Erstellen der Zeichenfläche mit paintComponent und Konstruktor Panel`r") {
        $inner = "<w:r><w:t xml:space='preserve'>Erstellen der Zeichenfläche mit </w:t></w:r>" +
                 "<w:proofErr w:type='spellStart'/>" +
                 "<w:r><w:t>paintComponent</w:t></w:r>" +
                 "<w:proofErr w:type='spellEnd'/>" +
                 "<w:r><w:t xml:space='preserve'> und Konstruktor Panel</w:t></w:r>"
        Set-ListParaRuns $p $inner
    }
}

# Append the new bullet point as the last item in the same bulleted list.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "Ändern der Farbe mit Buttons"
